$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Fri Aug 30 19:14:22 UTC 2024 with GitHub Actions

$ws.Range("D2").Value = "'58.777.80"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "'2.489.73"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'533.47"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'136.03"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.560"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("D9").Value = "'2.508.24"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").Value = "'0.100"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'0.159"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").Value = "'5.29"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "'0.346"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "'2.937.93"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "'23.00"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'58.751.10"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "'2.509.13"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "'11.00"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'4.24"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'324.28"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'63.29"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").Value = "'0.416"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "'7.56"
$ws.Range("E28").Value = "  -2.86%  "
$ws.Range("D29").Value = "'6.75"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "'0.0₃0769"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'1.76"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").Value = "'167.53"
$ws.Range("E32").Value = "  +2.86%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.998"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.13"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "'1.38"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "'18.43"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'4.08"
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("D38").Value = "'1.55"
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("D39").Value = "'36.63"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "'0.818"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").Value = "'3.60"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "'5.22"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").Value = "'277.20"
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("D44").Value = "'0.995"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'0.599"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "'10.84"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "'125.48"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'0.0923"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'0.0509"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "'0.0220"
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("D51").Value = "'17.39"
$ws.Range("E51").Value = "  +0.22%  "
